$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "conversión del día" summary text -------------
$hoja1 = $wb.Worksheets.Item("Hoja1")
$nuevoTexto = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 6.15 = 24969.25 pesos
✅ 24969.25 pesos = 6.13 = 970.4 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@
$hoja1.Range("A1").Value = $nuevoTexto

# --- tasas: refresh the scraped Binance / transfi rate figures ----------
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 162.6
$tasas.Range("O10").Value = 4060
$tasas.Range("N12").Value = 4075
$tasas.Range("O12").Value = 158.37
